# Auto-update predictions and index for 2025-10-23
# This script rewrites the fixtures sheet: updates existing match rows with
# final results, inserts newly-completed fixtures, and relocates the
# summary (index) formulas to the bottom of the now-larger table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Wed Oct 22"
$ws.Range("B2").Value = "Zenit St. Petersburg ✓ - FC Orenburg: 6:0"
$ws.Range("C2").Value = 2.65
$ws.Range("D2").Value = "Zenit St. Petersburg"
$ws.Range("E2").Value = 3.5
$ws.Range("F2").Value = "76%"
$ws.Range("G2").Value = "✓"
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = $false

# Row 3
$ws.Range("A3").Value = "Wed Oct 22"
$ws.Range("B3").Value = "Bayern Munich ✓ - Club Brugge KV: 4:0"
$ws.Range("C3").Value = 2.61
$ws.Range("D3").Value = "Bayern Munich"
$ws.Range("E3").Value = 3.5
$ws.Range("F3").Value = "75%"
$ws.Range("G3").Value = "✓"
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = $false

# Row 4
$ws.Range("A4").Value = "Wed Oct 22"
$ws.Range("B4").Value = "Galatasaray ✓ - FK Bodø/Glimt: 3:1"
$ws.Range("C4").Value = 3.41
$ws.Range("D4").Value = "Galatasaray"
$ws.Range("E4").Value = 4.5
$ws.Range("F4").Value = "74%"
$ws.Range("G4").Value = "✓"
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = $true

# Row 5
$ws.Range("A5").Value = "Wed Oct 22"
$ws.Range("B5").Value = "Odds BK - Lillestrøm SK ✓: 1:7"
$ws.Range("C5").Value = 2.51
$ws.Range("D5").Value = "Lillestrøm SK"
$ws.Range("E5").Value = 3.5
$ws.Range("F5").Value = "74%"
$ws.Range("G5").Value = "✓"
$ws.Range("H5").Value = 8
$ws.Range("I5").Value = $false

# Row 6
$ws.Range("A6").Value = "Wed Oct 22"
$ws.Range("B6").Value = "CF América ✓ - Puebla FC: 2:1"
$ws.Range("C6").Value = 2.22
$ws.Range("D6").Value = "CF América"
$ws.Range("E6").Value = 3.5
$ws.Range("F6").Value = "73%"
$ws.Range("G6").Value = "✓"
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = $true

# Row 7
$ws.Range("A7").Value = "Wed Oct 22"
$ws.Range("B7").Value = "Real Madrid ✓ - Juventus FC: 1:0"
$ws.Range("C7").Value = 2.94
$ws.Range("D7").Value = "Real Madrid"
$ws.Range("E7").Value = 3.5
$ws.Range("F7").Value = "72%"
$ws.Range("G7").Value = "✓"
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = $true

# Row 8
$ws.Range("A8").Value = "Wed Oct 22"
$ws.Range("B8").Value = "FC Goa - Al-Nassr FC ✓: 1:2"
$ws.Range("C8").Value = 2.12
$ws.Range("D8").Value = "Al-Nassr FC"
$ws.Range("E8").Value = 3.5
$ws.Range("F8").Value = "71%"
$ws.Range("G8").Value = "✓"
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = $true

# Row 9
$ws.Range("A9").Value = "Wed Oct 22"
$ws.Range("B9").Value = "FC Sardarapat ✓ - FC Andranik: 1:0"
$ws.Range("C9").Value = 3.78
$ws.Range("D9").Value = "FC Sardarapat"
$ws.Range("E9").Value = 4.5
$ws.Range("F9").Value = "71%"
$ws.Range("G9").Value = "✓"
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = $true

# Row 10
$ws.Range("A10").Value = "Wed Oct 22"
$ws.Range("B10").Value = "Aris Limassol ✓ - Omonia 29is Maiou: 2:1"
$ws.Range("C10").Value = 3.05
$ws.Range("D10").Value = "Aris Limassol"
$ws.Range("E10").Value = 4.5
$ws.Range("F10").Value = "70%"
$ws.Range("G10").Value = "✓"
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = $true

# Row 11
$ws.Range("A11").Value = "Wed Oct 22"
$ws.Range("B11").Value = "Chelsea FC ✓ - Ajax Amsterdam: 5:1"
$ws.Range("C11").Value = 2.18
$ws.Range("D11").Value = "Chelsea FC"
$ws.Range("E11").Value = 3.5
$ws.Range("F11").Value = "70%"
$ws.Range("G11").Value = "✓"
$ws.Range("H11").Value = 6
$ws.Range("I11").Value = $false

# Row 12
$ws.Range("A12").Value = "Wed Oct 22"
$ws.Range("B12").Value = "Kongsvinger IL ✓ - IL Hødd: 3:0"
$ws.Range("C12").Value = 3.89
$ws.Range("D12").Value = "Kongsvinger IL"
$ws.Range("E12").Value = 4.5
$ws.Range("F12").Value = "66%"
$ws.Range("G12").Value = "✓"
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = $true

# Row 13
$ws.Range("A13").Value = "Wed Oct 22"
$ws.Range("B13").Value = "FCI Levadia X - Paide Linnameeskond: 1:2"
$ws.Range("C13").Value = 1.92
$ws.Range("D13").Value = "FCI Levadia"
$ws.Range("E13").Value = 2.5
$ws.Range("F13").Value = "64%"
$ws.Range("G13").Value = "X"
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = $false

# Row 14
$ws.Range("A14").Value = "Wed Oct 22"
$ws.Range("B14").Value = "Kalev Tallinn - Kalju FC ✓: 0:2"
$ws.Range("C14").Value = 3.76
$ws.Range("D14").Value = "Kalju FC"
$ws.Range("E14").Value = 4.5
$ws.Range("F14").Value = "63%"
$ws.Range("G14").Value = "✓"
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = $true

# Row 15
$ws.Range("A15").Value = "Wed Oct 22"
$ws.Range("B15").Value = "Aalesunds FK ✓ - Lyn 1896 FK: 2:1"
$ws.Range("C15").Value = 1.84
$ws.Range("D15").Value = "Aalesunds FK"
$ws.Range("E15").Value = 2.5
$ws.Range("F15").Value = "60%"
$ws.Range("G15").Value = "✓"
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = $false

# Row 16
$ws.Range("A16").Value = "Wed Oct 22"
$ws.Range("B16").Value = "FC Santa Coloma ✓ - Penya Encarnada d'Andorra: 1:0"
$ws.Range("C16").Value = 2.43
$ws.Range("D16").Value = "FC Santa Coloma"
$ws.Range("E16").Value = 3.5
$ws.Range("F16").Value = "60%"
$ws.Range("G16").Value = "✓"
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = $true

# Row 17
$ws.Range("A17").Value = "Wed Oct 22"
$ws.Range("B17").Value = "Nagaworld FC - Phnom Penh Crown : 12:00"
$ws.Range("C17").Value = 2.21
$ws.Range("D17").Value = "Phnom Penh Crown"
$ws.Range("E17").Value = 3.5
$ws.Range("F17").Value = "60%"
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = 12
$ws.Range("I17").Value = $false

# Row 18
$ws.Range("A18").Value = "Wed Oct 22"
$ws.Range("B18").Value = "JS Kabylie  - MC El Bayadh: 23:00"
$ws.Range("C18").Value = 0.88
$ws.Range("D18").Value = "JS Kabylie"
$ws.Range("E18").Value = 1.5
$ws.Range("F18").Value = "58%"
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = 23
$ws.Range("I18").Value = $false

# Row 19
$ws.Range("A19").Value = "Wed Oct 22"
$ws.Range("B19").Value = "Stabæk Fotball ✓ - Åsane Fotball: 2:0"
$ws.Range("C19").Value = 2.73
$ws.Range("D19").Value = "Stabæk Fotball"
$ws.Range("E19").Value = 3.5
$ws.Range("F19").Value = "56%"
$ws.Range("G19").Value = "✓"
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = $true

# Row 20
$ws.Range("A20").Value = "Wed Oct 22"
$ws.Range("B20").Value = "MC Algiers  - JS Saoura: 23:00"
$ws.Range("C20").Value = 1.07
$ws.Range("D20").Value = "MC Algiers"
$ws.Range("E20").Value = 2.5
$ws.Range("F20").Value = "55%"
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = 23
$ws.Range("I20").Value = $false

# Row 21
$ws.Range("A21").Value = "Wed Oct 22"
$ws.Range("B21").Value = "Deportivo La Guaira  - Zamora FC: 1:1"
$ws.Range("C21").Value = 1.37
$ws.Range("D21").Value = "Deportivo La Guaira"
$ws.Range("E21").Value = 2.5
$ws.Range("F21").Value = "53%"
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = $true

# Row 22
$ws.Range("A22").Value = "Wed Oct 22"
$ws.Range("B22").Value = "MC Algiers ✓ - Paradou AC: 2:1"
$ws.Range("C22").Value = 1.53
$ws.Range("D22").Value = "MC Algiers"
$ws.Range("E22").Value = 2.5
$ws.Range("F22").Value = "52%"
$ws.Range("G22").Value = "✓"
$ws.Range("H22").Value = 3
$ws.Range("I22").Value = $false


# The summary/index block used to live at K12:L14; the table has grown to
# 21 match rows (rows 2-22), so relocate it to K23:L25 and drop the old
# cells so they don't linger as stale leftovers.
$ws.Range("K12:L14").ClearContents()

$ws.Range("K23").Formula = "=COUNTIF(I:I,TRUE)"
$ws.Range("L23").Formula = "=(K23/K25)*100"
$ws.Range("K24").Formula = "=COUNTIF(I:I,FALSE)"
$ws.Range("K25").Formula = "=K23+K24"
